$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 419
$ws.Range("A419").Value = "Amon"
$ws.Range("B419").Value = "longitude latitude time"
$ws.Range("C419").Value = "fco2antt"
$ws.Range("D419").Value = "Carbon Mass Flux into Atmosphere Due to All Anthropogenic Emissions of CO2"
$ws.Range("E419").Value = "Not available in the AOGCM, but will be added by Tommi  in the ESM in TM5 with its cmor name."
$ws.Range("F419").Value = "Tommi Bergman"
$ws.Range("H419").Value = "This is requested only for the emission-driven coupled carbon climate model runs.  Does not include natural fire sources but, includes all anthropogenic sources, including fossil fuel use, cement production, agricultural burning, and sources associated with anthropogenic land use change excluding forest regrowth."

# Row 420
$ws.Range("A420").Value = "Amon"
$ws.Range("B420").Value = "longitude latitude time"
$ws.Range("C420").Value = "fco2fos"
$ws.Range("D420").Value = "Carbon Mass Flux into Atmosphere Due to Fossil Fuel Emissions of CO2"
$ws.Range("E420").Value = "Not available in the AOGCM, but will be added by Tommi  in the ESM in TM5 with its cmor name."
$ws.Range("F420").Value = "Tommi Bergman"
$ws.Range("H420").Value = "This is the prescribed anthropogenic CO2 flux from fossil fuel use, including cement production, and flaring (but not from land-use changes, agricultural burning, forest regrowth, etc.)"

# Row 421
$ws.Range("A421").Value = "Amon"
$ws.Range("B421").Value = "longitude latitude time"
$ws.Range("C421").Value = "fco2nat"
$ws.Range("D421").Value = "Surface Carbon Mass Flux into the Atmosphere Due to Natural Sources"
$ws.Range("E421").Value = "Not available in the AOGCM, but will be added by Tommi  in the ESM in TM5 with its cmor name."
$ws.Range("F421").Value = "Tommi Bergman"
$ws.Range("H421").Value = "This is what the atmosphere sees (on its own grid).  This field should be equivalent to the combined natural fluxes of carbon  that account for natural exchanges between the atmosphere and land (nep) or ocean (fgco2) reservoirs."

# Row 422
$ws.Range("A422").Value = "Amon"
$ws.Range("B422").Value = "longitude latitude plev19 time2"
$ws.Range("C422").Value = "o3Clim"
$ws.Range("D422").Value = "Mole Fraction of O3"
$ws.Range("E422").Value = "Not available in the AOGCM, but will be added by Tommi  in the ESM in TM5 with its cmor name."
$ws.Range("F422").Value = "Tommi Bergman"

# Row 423
$ws.Range("A423").Value = "Amon"
$ws.Range("B423").Value = "longitude latitude plev19 time2"
$ws.Range("C423").Value = "co2Clim"
$ws.Range("D423").Value = "Mole Fraction of CO2"
$ws.Range("E423").Value = "Not available in the AOGCM, but will be added by Tommi  in the ESM in TM5 with its cmor name."
$ws.Range("F423").Value = "Tommi Bergman"

# Row 424
$ws.Range("A424").Value = "Amon"
$ws.Range("B424").Value = "time"
$ws.Range("C424").Value = "co2mass"
$ws.Range("D424").Value = "Total Atmospheric Mass of CO2"
$ws.Range("E424").Value = "Not available in the AOGCM, but will be added by Tommi  in the ESM in TM5 with its cmor name."
$ws.Range("F424").Value = "Tommi Bergman"
$ws.Range("H424").Value = "Total atmospheric mass of Carbon Dioxide"

# Row 425
$ws.Range("A425").Value = "Amon"
$ws.Range("B425").Value = "time2"
$ws.Range("C425").Value = "co2massClim"
$ws.Range("D425").Value = "Total Atmospheric Mass of CO2"
$ws.Range("E425").Value = "Not available in the AOGCM, but will be added by Tommi  in the ESM in TM5 with its cmor name."
$ws.Range("F425").Value = "Tommi Bergman"
$ws.Range("H425").Value = "Total atmospheric mass of Carbon Dioxide"

# Row 426
$ws.Range("A426").Value = "Amon"
$ws.Range("B426").Value = "longitude latitude plev19 time2"
$ws.Range("C426").Value = "ch4Clim"
$ws.Range("D426").Value = "Mole Fraction of CH4"
$ws.Range("E426").Value = "Not available in the AOGCM, but will be added by Tommi  in the ESM in TM5 with its cmor name."
$ws.Range("F426").Value = "Tommi Bergman"

# Row 427
$ws.Range("A427").Value = "Amon"
$ws.Range("B427").Value = "time"
$ws.Range("C427").Value = "ch4global"
$ws.Range("D427").Value = "Global Mean Mole Fraction of CH4"
$ws.Range("E427").Value = "Not available in the AOGCM, but will be added by Tommi  in the ESM in TM5 with its cmor name."
$ws.Range("F427").Value = "Tommi Bergman"
$ws.Range("H427").Value = "Global Mean Mole Fraction of CH4"

# Row 428
$ws.Range("A428").Value = "Amon"
$ws.Range("B428").Value = "time2"
$ws.Range("C428").Value = "ch4globalClim"
$ws.Range("D428").Value = "Global Mean Mole Fraction of CH4"
$ws.Range("E428").Value = "Not available in the AOGCM, but will be added by Tommi  in the ESM in TM5 with its cmor name."
$ws.Range("F428").Value = "Tommi Bergman"
$ws.Range("H428").Value = "Global Mean Mole Fraction of CH4"

# Row 429
$ws.Range("A429").Value = "Amon"
$ws.Range("B429").Value = "longitude latitude plev19 time2"
$ws.Range("C429").Value = "n2oClim"
$ws.Range("D429").Value = "Mole Fraction of N2O"
$ws.Range("E429").Value = "Not available in the AOGCM, but will be added by Tommi  in the ESM in TM5 with its cmor name."
$ws.Range("F429").Value = "Tommi Bergman"

# Row 430
$ws.Range("A430").Value = "Amon"
$ws.Range("B430").Value = "time"
$ws.Range("C430").Value = "n2oglobal"
$ws.Range("D430").Value = "Global Mean Mole Fraction of N2O"
$ws.Range("E430").Value = "Not available in the AOGCM, but will be added by Tommi  in the ESM in TM5 with its cmor name."
$ws.Range("F430").Value = "Tommi Bergman"
$ws.Range("H430").Value = "Global mean Nitrous Oxide (N2O)"

# Row 431
$ws.Range("A431").Value = "Amon"
$ws.Range("B431").Value = "time2"
$ws.Range("C431").Value = "n2oglobalClim"
$ws.Range("D431").Value = "Global Mean Mole Fraction of N2O"
$ws.Range("E431").Value = "Not available in the AOGCM, but will be added by Tommi  in the ESM in TM5 with its cmor name."
$ws.Range("F431").Value = "Tommi Bergman"
$ws.Range("H431").Value = "Global mean Nitrous Oxide (N2O)"

$ws.Range("F419:F431").Select()

# Extend the sheet's recorded dimension out to the full column/row extent
# (mirrors the wider-than-data-range dimension the workbook ends up with
# after the interactive edit session in the source diff) without
# disturbing any visible cell content or styling.
$ws.Cells.Item(1048576, 8).NumberFormat = "General"
